$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.057.39'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.227.63'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '294.30'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '87.34'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.64%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.472'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '30.19'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0780'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.78%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.47'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.570.80'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.83%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '13.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.236.59'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.728'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '39.959.42'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0886'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.31'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +5.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.81'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.52'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.51%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.14'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.28%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.00%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.94%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.73'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.10'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -5.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.20'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '156.21'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.86%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '31.72'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -7.32%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.06%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.19%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.91'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +6.88%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.69%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '15.73'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0971'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -5.11%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.68'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.123.79'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.22%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.13'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -4.60%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '18.11'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +10.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.79'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.65'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.434.34'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.14%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.23%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '88.73'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.57%  '
